$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 155; this shifts existing rows 155-242 down to 156-243
$ws.Rows.Item(155).Insert()

# Populate the new row 155 with data (copy structural/text columns from the row below,
# which now holds what used to be row 155, and set the new values per the edit)
$ws.Cells.Item(155, 1).Value = 8
$ws.Cells.Item(155, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44830
$ws.Cells.Item(155, 5).Value = 4
$ws.Cells.Item(155, 6).Value = 100112037
$ws.Cells.Item(155, 7).Value = "Cebollín"
$ws.Cells.Item(155, 8).Value = "Sin especificar"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 1200
$ws.Cells.Item(155, 11).Value = 1400
$ws.Cells.Item(155, 12).Value = 1600
$ws.Cells.Item(155, 13).Value = 1500
$ws.Cells.Item(155, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(155, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(155, 16).Value = 250
$ws.Cells.Item(155, 17).Value = 6
$ws.Cells.Item(155, 18).Value = "Hortaliza"

# Ensure the date column keeps the same number format style as the rest of column D
$ws.Cells.Item(155, 4).NumberFormat = $ws.Cells.Item(156, 4).NumberFormat
